$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 6, shifting existing rows 6-15 down to 7-16
$ws.Rows.Item(6).Insert()

# Fill in the new row's data
$ws.Range("A6").Value = "Age"
$ws.Range("B6").Value = "Age of cars in years"
$ws.Range("C6").Value = "numeric"

# Move the selection to A7, matching the final state of the workbook
$ws.Range("A7").Select()
